# Add analysed results for medium density - Reference
# Fills in the "Reference" scenario rows (rows 3-12) of the "Medium traffic
# density" worksheet, and moves the active sheet/selection the way the
# author left the workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Low traffic densit
$ws2 = $wb.Worksheets.Item(2)   # Medium traffic density

# ---------------------------------------------------------------------
# Raw data for the "Reference" rows (A3:B12) of the Medium traffic table
# ---------------------------------------------------------------------
$ws2.Range("C3").Value = 230
$ws2.Range("D3").Value = 182
$ws2.Range("E3").Value = 3728.2939018162401
$ws2.Range("F3").Value = 313.22747252747303
$ws2.Range("H3").Value = 0.943585304942688
$ws2.Range("I3").Value = 59
$ws2.Range("J3").Value = 6.6779661016949197
$ws2.Range("K3").Value = 31
$ws2.Range("L3").Value = 5

$ws2.Range("C4").Value = 230
$ws2.Range("D4").Value = 182
$ws2.Range("E4").Value = 3645.57271517038
$ws2.Range("F4").Value = 310.99862637362702
$ws2.Range("H4").Value = 0.98136119459873705
$ws2.Range("I4").Value = 57
$ws2.Range("J4").Value = 17.315789473684202
$ws2.Range("K4").Value = 32
$ws2.Range("L4").Value = 12

$ws2.Range("C5").Value = 227
$ws2.Range("D5").Value = 180
$ws2.Range("E5").Value = 3468.5506895009398
$ws2.Range("F5").Value = 292.17750000000001
$ws2.Range("H5").Value = 0.90648305124000605
$ws2.Range("I5").Value = 46
$ws2.Range("J5").Value = 23.630434782608699
$ws2.Range("K5").Value = 33
$ws2.Range("L5").Value = 4

$ws2.Range("C6").Value = 225
$ws2.Range("D6").Value = 181
$ws2.Range("E6").Value = 3556.5964387519298
$ws2.Range("F6").Value = 296.99861878452998
$ws2.Range("H6").Value = 0.87445410350742303
$ws2.Range("I6").Value = 55
$ws2.Range("J6").Value = 8.2909090899999995
$ws2.Range("K6").Value = 34
$ws2.Range("L6").Value = 9

$ws2.Range("C7").Value = 236
$ws2.Range("D7").Value = 183
$ws2.Range("E7").Value = 3541.6060530570498
$ws2.Range("F7").Value = 291.38852459016402
$ws2.Range("H7").Value = 0.99115935526818999
$ws2.Range("I7").Value = 58
$ws2.Range("J7").Value = 14.689655172413801
$ws2.Range("K7").Value = 38
$ws2.Range("L7").Value = 11

$ws2.Range("C8").Value = 221
$ws2.Range("D8").Value = 178
$ws2.Range("E8").Value = 3551.55124933562
$ws2.Range("F8").Value = 313.60842696629197
$ws2.Range("H8").Value = 0.90045219843297897
$ws2.Range("I8").Value = 48
$ws2.Range("J8").Value = 9.4166666669999994
$ws2.Range("K8").Value = 32
$ws2.Range("L8").Value = 7

$ws2.Range("C9").Value = 225
$ws2.Range("D9").Value = 183
$ws2.Range("E9").Value = 3519.35645950656
$ws2.Range("F9").Value = 298.412568306011
$ws2.Range("H9").Value = 0.88235792606951802
$ws2.Range("I9").Value = 62
$ws2.Range("J9").Value = 9.5806451612903203
$ws2.Range("K9").Value = 38
$ws2.Range("L9").Value = 10

$ws2.Range("C10").Value = 230
$ws2.Range("D10").Value = 180
$ws2.Range("E10").Value = 3565.58677805467
$ws2.Range("F10").Value = 300.58999999999997
$ws2.Range("H10").Value = 0.91530053560013802
$ws2.Range("I10").Value = 58
$ws2.Range("J10").Value = 7.2241379310344804
$ws2.Range("K10").Value = 38
$ws2.Range("L10").Value = 8

$ws2.Range("C11").Value = 234
$ws2.Range("D11").Value = 186
$ws2.Range("E11").Value = 3635.7007822218802
$ws2.Range("F11").Value = 318.61881720430102
$ws2.Range("H11").Value = 1.02991938215926
$ws2.Range("I11").Value = 60
$ws2.Range("J11").Value = 10.683333299999999
$ws2.Range("K11").Value = 36
$ws2.Range("L11").Value = 15

$ws2.Range("C12").Value = 225
$ws2.Range("D12").Value = 182
$ws2.Range("E12").Value = 3436.7056426013201
$ws2.Range("F12").Value = 301.74725274725301
$ws2.Range("H12").Value = 0.94134338096429204
$ws2.Range("I12").Value = 64
$ws2.Range("J12").Value = 10.0625
$ws2.Range("K12").Value = 41
$ws2.Range("L12").Value = 8

# ---------------------------------------------------------------------
# Force the dependent formulas (G3:G12 ratio column, and the AVERAGE /
# STDEV.P summary rows 13 & 14) to recompute against the new inputs, so
# that the cached <v> values stored in the file are refreshed.
# ---------------------------------------------------------------------
$null = $ws2.Range("G3:L14").Value2

# ---------------------------------------------------------------------
# Update sheet selections / the active tab: the workbook now opens on
# the Medium traffic density sheet, with both sheets' selections left
# at A15.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A15").Select()

$ws2.Activate()
$ws2.Range("A15").Select()
